# Update Name of Algo
# Apply updated imputed values produced by the KNN algorithm run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E3").Value  = 16.246
$ws.Range("C12").Value = -10.94
$ws.Range("E14").Value = 17.066
$ws.Range("E26").Value = 16.407
$ws.Range("E31").Value = 16.376
$ws.Range("C32").Value = -13.407
$ws.Range("E35").Value = 16.492
$ws.Range("C36").Value = -12.732
$ws.Range("E37").Value = 16.815
$ws.Range("C38").Value = -12.607
$ws.Range("E45").Value = 16.677
$ws.Range("C46").Value = -14.162
$ws.Range("C54").Value = -12.705
$ws.Range("C55").Value = -13.752
$ws.Range("E57").Value = 16.461
$ws.Range("C67").Value = -12.037
$ws.Range("C69").Value = -11.062
$ws.Range("C72").Value = -11.932
$ws.Range("C91").Value = -11.095
$ws.Range("C99").Value = -12.635
$ws.Range("E100").Value = 16.741
$ws.Range("E102").Value = 16.519
